# Apply sample-data corrections to Sheet1:
#  - Fix typo "Eric Carman" -> "Eric Cartman"
#  - Insert a new row for Peggy Hill (Hank Hill's spouse)
#  - Append a new row for Marge Simpson (Homer's spouse)
#  - Update the active selection to a single cell

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the typo'd name in A3
$ws.Range("A3").Value = "Eric Cartman"

# Insert a new row at position 4 (shifts the old rows 4 and 5 down to 5 and 6)
$ws.Rows("4").Insert()

# Populate the newly inserted row 4 with Peggy Hill's data
$ws.Range("A4").Value = "Peggy Hill"
$ws.Range("B4").Value = 38
$ws.Range("C4").Value = 66
$ws.Range("D4").Value = 125
$ws.Range("E4").Value = "Hank Hill"

# Append a new row 7 with Marge Simpson's data
$ws.Range("A7").Value = "Marge Simpson"
$ws.Range("B7").Value = 43
$ws.Range("C7").Value = 78
$ws.Range("D7").Value = 135
$ws.Range("E7").Value = "Homer J. Simpson"

# Update the selected cell to match the new layout
$ws.Range("E5").Select()
